# Lisa_Fleischmann_Wink_Tour_Rev3.pptx — "Rev to Week 4 and 5 remove numbered steps"
#
# 1) Update the cached "datetimeFigureOut" field text (2/20/2020 -> 4/3/2020)
#    on every Slide Master / Custom Layout Date placeholder.
# 2) Remove the numbered red-circle callout ovals ("1".."5") from slides
#    3, 4, 5 (x2) and 6.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text refresh across every master + layout.
# ---------------------------------------------------------------------------
$oldDate = "2/20/2020"
$newDate = "4/3/2020"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    Update-DatePlaceholders $master.Shapes

    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        $layout = $master.CustomLayouts.Item($li)
        Update-DatePlaceholders $layout.Shapes
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the numbered oval callouts.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item("Oval 9").Delete()          # "1"

$s4 = $p.Slides.Item(4)
$s4.Shapes.Item("Oval 7").Delete()          # "2"

$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("Oval 8").Delete()          # "4" (delete higher index first)
$s5.Shapes.Item("Oval 6").Delete()          # "3"

$s6 = $p.Slides.Item(6)
$s6.Shapes.Item("Oval 6").Delete()          # "5"
